$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange
$para = $tr.Paragraphs(3)

# Paragraph currently reads:
#   "d" (run: lang=en-US)  + "er gleiche Bearbeitungsauftrag mehrmals mit verschiedenen Datenkonstellationen getestet werden" (run: lang=de-DE)
# Target:
#   "D" (run: lang=en-US) + "er " (run: lang=de-DE) + "gleiche Bearbeitungsauftrag mehrmals mit verschiedenen Datenkonstellationen getestet werden" (run: lang=de-DE)

# 1) Fix capitalisation of the first run ("d" -> "D"), keeping its formatting.
$firstChar = $para.Characters(1, 1)
$firstChar.Text = "D"

# 2) Split the second run into "er " and "gleiche ..." by re-asserting the text of the
#    leading "er " portion. Re-assigning the same text to a sub-range of a run forces the
#    host to materialise it as its own run while leaving the original run formatting intact.
$erPortion = $para.Characters(2, 3)
$erPortion.Text = "er "
